$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null) {
        $trimmed = $val -replace '[\. ]+$', ''
        if ($trimmed -ne $val) {
            $cell.Value = $trimmed
        }
    }
}
